$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 4.9842303083526174
$ws.Range("C2").Value = 0.90999355879820343
$ws.Range("D2").Value = 0.68550585053288438
$ws.Range("E2").Value = 0.10245382994170703

$ws.Range("B3").Value = 4.1911312703011223
$ws.Range("C3").Value = 5.8905749469724293
$ws.Range("D3").Value = 5.3759646401996122
$ws.Range("E3").Value = -1.9344935784110526

$ws.Range("B1:E3").Select()
